$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-11-18 Monday" "2024-11-19 Tuesday"

Replace-Text "415×2=830" "720×4=2880"
Replace-Text "489×4=1956" "921×2=1842"
Replace-Text "520×3=1560" "113×9=1017"
Replace-Text "882×3=2646" "563×5=2815"
Replace-Text "251×6=1506" "930×3=2790"

Replace-Text "613×2=1226" "469×7=3283"
Replace-Text "869×9=7821" "747×5=3735"
Replace-Text "811×7=5677" "280×5=1400"
Replace-Text "135×3=405" "523×5=2615"
Replace-Text "527×6=3162" "484×5=2420"

Replace-Text "772×8=6176" "756×7=5292"
Replace-Text "448×8=3584" "331×6=1986"
Replace-Text "146×5=730" "536×3=1608"
Replace-Text "909×2=1818" "730×6=4380"
Replace-Text "624×4=2496" "910×9=8190"

Replace-Text "434×5=2170" "406×5=2030"
Replace-Text "633×2=1266" "400×5=2000"
Replace-Text "253×4=1012" "172×9=1548"
Replace-Text "955×8=7640" "122×4=488"
Replace-Text "176×3=528" "389×3=1167"

Replace-Text "342×2=684" "495×7=3465"
Replace-Text "173×2=346" "716×9=6444"
Replace-Text "675×7=4725" "696×3=2088"
Replace-Text "582×6=3492" "596×3=1788"
Replace-Text "962×5=4810" "911×3=2733"
